$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TypeMapping")

# 1. Rename the "Millimeter" unit to "Centimeter" (A26) and change its JType from Integer to Float (B26)
$ws.Range("A26").Value = "Centimeter"
$ws.Range("B26").Value = "Float"

# 2. Add a new PIM/JType mapping row for IpCimMaszk -> String
$ws.Range("A30").Value = "IpCimMaszk"
$ws.Range("B30").Value = "String"

# 3. Move the active selection/cursor
$ws.Range("K18").Select()
